# Apply the 2023-10-22 cryptos-list refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.931.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0608"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +18.12%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.966.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("E23").Value = "  +2.39%  "
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.422.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  +4.28%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.828"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.775.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  +8.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.57%  "
